$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -11.4974
$ws.Range("B7").Value = 4.737599999999993
$ws.Range("A8").Value = -22.40650000000002
$ws.Range("A10").Value = -21.73999999999998
$ws.Range("E10").Value = 16.52979999999999
$ws.Range("A12").Value = -21.56470000000001
$ws.Range("E12").Value = 17.95220000000002
$ws.Range("E13").Value = 16.4509
$ws.Range("E14").Value = 16.8527
$ws.Range("B15").Value = 4.692099999999994
$ws.Range("A18").Value = -21.80629999999999
$ws.Range("B18").Value = 6.314400000000003
$ws.Range("C18").Value = -12.267
$ws.Range("C19").Value = -11.44100000000001
$ws.Range("B20").Value = 9.289400000000002
$ws.Range("C27").Value = -13.25089999999999
$ws.Range("B29").Value = 5.149800000000003
$ws.Range("E29").Value = 17.32650000000001
$ws.Range("B30").Value = 4.510500000000001
$ws.Range("B31").Value = 4.987399999999999
$ws.Range("C31").Value = -13.33159999999999
$ws.Range("E32").Value = 16.75849999999998
$ws.Range("E35").Value = 16.22530000000001
$ws.Range("A37").Value = -19.71189999999999
$ws.Range("C38").Value = -13.1849
$ws.Range("B40").Value = 9.442799999999993
$ws.Range("C42").Value = -12.03190000000001
$ws.Range("E43").Value = 17.18430000000001
$ws.Range("C44").Value = -14.00869999999999
$ws.Range("C47").Value = -12.4626
$ws.Range("E48").Value = 17.33510000000001
$ws.Range("E49").Value = 15.6477
$ws.Range("B50").Value = 5.26
$ws.Range("E50").Value = 16.33929999999999
$ws.Range("A55").Value = -21.7392
$ws.Range("E56").Value = 16.3984
$ws.Range("C58").Value = -11.67429999999998
$ws.Range("C65").Value = -12.22879999999999
$ws.Range("A68").Value = -21.4683
$ws.Range("B68").Value = 4.566799999999995
$ws.Range("E69").Value = 17.33040000000003
$ws.Range("C73").Value = -12.34990000000001
$ws.Range("B76").Value = 6.396699999999997
$ws.Range("A77").Value = -20.11029999999998
$ws.Range("A78").Value = -20.14999999999998
$ws.Range("A81").Value = -22.10680000000002
$ws.Range("E81").Value = 16.81809999999999
$ws.Range("A82").Value = -21.6882
$ws.Range("B87").Value = 4.524699999999993
$ws.Range("B88").Value = 4.647399999999997
$ws.Range("C90").Value = -13.0003
$ws.Range("E92").Value = 18.22600000000001
$ws.Range("C94").Value = -10.0169
$ws.Range("C95").Value = -12.5011
$ws.Range("B96").Value = 5.175900000000007
$ws.Range("B98").Value = 7.407299999999998
$ws.Range("B101").Value = 9.4413
$ws.Range("C101").Value = -12.7327
$ws.Range("B102").Value = 8.665000000000006
